$wb = $excel.ActiveWorkbook

# Remember which sheet/tab was active so we can restore it at the end -
# selecting a range on another sheet necessarily activates that sheet.
$origActiveSheet = $wb.ActiveSheet

$ws = $wb.Worksheets.Item("Service Contacts")

# The "delivery_organisation_path" column (R) needs to move so that it sits
# immediately before "practitioner_key" (currently D), i.e. become the new
# column D. Cut the whole column and insert it (with its cell formatting)
# in front of column D - this shifts the existing D:Q columns right into
# E:R, exactly like Excel's "Cut" + "Insert Cut Cells" on whole columns.
$ws.Columns("R").Cut() | Out-Null
$ws.Columns("D").Insert() | Out-Null

# After an "Insert Cut Cells" operation Excel leaves the newly inserted
# column selected.
$ws.Range("D1:D1048576").Select() | Out-Null

# Restore the original active sheet/tab.
$origActiveSheet.Activate() | Out-Null
